# [update] : Player Stat 레벨링 저장 /불러오기 구현
# Adds two new columns (AA: SoulIDList, AB: ItemIDList) to the UserDB sheet,
# mirroring the existing header/data formatting, plus sample data in row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy existing column formatting onto the new AA/AB columns --------
# Z1 carries the header style (bold, centered); Z2/Z3/Z4 carry the plain
# data style (vertical-center only) used across the whole data block.
$ws.Range("Z1").Copy()
$ws.Range("AA1:AB1").PasteSpecial(-4122)

$ws.Range("Z2").Copy()
$ws.Range("AA2:AB3").PasteSpecial(-4122)
$ws.Range("AB4").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Header row (row 1) --------------------------------------------------
$ws.Range("AA1").Value = "SoulIDList"
$ws.Range("AB1").Value = "ItemIDList"

# --- Type row (row 2) -----------------------------------------------------
$ws.Range("AA2").Value = "List<int>"
$ws.Range("AB2").Value = "List<int>"

# --- Description row (row 3), mirrors header key names ----------------
$ws.Range("AA3").Value = "SoulIDList"
$ws.Range("AB3").Value = "ItemIDList"

# --- Sample data row (row 4) ----------------------------------------------
$ws.Range("AA4").Value = 110001100111002
$ws.Range("AA4").NumberFormat = "#,##0"
$ws.Range("AB4").Value = 1000

# --- Column widths (best-fit sized to the new header text) ---------------
$ws.Columns.Item(27).ColumnWidth = 18.14
$ws.Columns.Item(28).ColumnWidth = 10.43

# --- Sheet view: scroll/zoom to show the new columns, matching the author's
#     final on-screen state when they saved the workbook -------------------
$ws.Application.ActiveWindow.Zoom = 115
[void]$ws.Range("AB10").Select()
